$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText($cellRef, $text) {
    $c = $ws.Range($cellRef)
    $origStyle = $c.Style
    $c.NumberFormat = "@"
    $c.Value = $text
    $c.Style = $origStyle
}

Set-CellText 'D2' '58.239.69'
Set-CellText 'E2' '  -3.56%  '
Set-CellText 'D3' '3.137.98'
Set-CellText 'E3' '  -5.13%  '
Set-CellText 'E4' '  +0.02%  '
Set-CellText 'D5' '523.87'
Set-CellText 'E5' '  -6.16%  '
Set-CellText 'D6' '134.99'
Set-CellText 'E6' '  -5.27%  '
Set-CellText 'E7' '  +0.04%  '
Set-CellText 'D8' '3.138.32'
Set-CellText 'E8' '  -5.12%  '
Set-CellText 'E9' '  -5.16%  '
Set-CellText 'D10' '7.26'
Set-CellText 'E10' '  -7.53%  '
Set-CellText 'D11' '0.108'
Set-CellText 'E11' '  -8.87%  '
Set-CellText 'E12' '  -6.48%  '
Set-CellText 'D13' '3.676.27'
Set-CellText 'E13' '  -5.12%  '
Set-CellText 'D15' '25.60'
Set-CellText 'E15' '  -4.72%  '
Set-CellText 'D16' '3.137.38'
Set-CellText 'E16' '  -4.98%  '
Set-CellText 'D17' '58.226.66'
Set-CellText 'E17' '  -3.59%  '
Set-CellText 'D18' '0.0000152'
Set-CellText 'E18' '  -7.79%  '
Set-CellText 'D19' '5.79'
Set-CellText 'E19' '  -5.28%  '
Set-CellText 'D20' '13.11'
Set-CellText 'E20' '  -7.25%  '
Set-CellText 'D21' '7.95'
Set-CellText 'E21' '  -8.33%  '
Set-CellText 'D22' '344.31'
Set-CellText 'E22' '  -8.24%  '
Set-CellText 'E23' '  -0.04%  '
Set-CellText 'D24' '68.65'
Set-CellText 'E24' '  -7.91%  '
Set-CellText 'D25' '0.507'
Set-CellText 'E25' '  -5.60%  '
Set-CellText 'D26' '3.267.52'
Set-CellText 'E26' '  -5.18%  '
Set-CellText 'E27' '  -1.92%  '
Set-CellText 'D28' '0.0₃0955'
Set-CellText 'E28' '  -7.06%  '
Set-CellText 'E29' '  +0.59%  '
Set-CellText 'D30' '6.80'
Set-CellText 'E30' '  -5.27%  '
Set-CellText 'D31' '0.999'
Set-CellText 'E31' '  -0.01%  '
Set-CellText 'E32' '  -8.88%  '
Set-CellText 'D33' '6.86'
Set-CellText 'E33' '  -9.55%  '
Set-CellText 'D34' '21.51'
Set-CellText 'E34' '  -5.03%  '
Set-CellText 'E35' '  -1.71%  '
Set-CellText 'D36' '4.81'
Set-CellText 'E36' '  -6.10%  '
Set-CellText 'D37' '157.03'
Set-CellText 'E37' '  -5.77%  '
Set-CellText 'D38' '6.22'
Set-CellText 'E38' '  -7.22%  '
Set-CellText 'E39' '  -10.22%  '
Set-CellText 'D40' '0.0692'
Set-CellText 'E40' '  -5.00%  '
Set-CellText 'D41' '3.170.20'
Set-CellText 'E41' '  -5.07%  '
Set-CellText 'B42' 'EnergySwap'
Set-CellText 'C42' 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-CellText 'D42' '24.40'
Set-CellText 'E42' '  -9.08%  '
Set-CellText 'B43' 'OKB'
Set-CellText 'C43' 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-CellText 'D43' '40.43'
Set-CellText 'E43' '  -3.57%  '
Set-CellText 'D44' '0.692'
Set-CellText 'E44' '  -7.91%  '
Set-CellText 'E45' '  -2.46%  '
Set-CellText 'E46' '  -5.85%  '
Set-CellText 'E47' '  -0.02%  '
Set-CellText 'E48' '  -8.74%  '
Set-CellText 'D49' '2.258.75'
Set-CellText 'E49' '  -4.26%  '
Set-CellText 'D50' '6.19'
Set-CellText 'E50' '  -3.51%  '
Set-CellText 'D51' '20.55'
Set-CellText 'E51' '  -3.63%  '
